# Update cryptocurrency price/volume figures per the latest data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.947.28'
$ws.Range('E2').Value = '  +0.46%  '
$ws.Range('D3').Value = '2.342.64'
$ws.Range('E3').Value = '  +0.58%  '
$ws.Range('E4').Value = '  -0.45%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '540.93'
$c.ClearFormats()
$ws.Range('E5').Value = '  -0.14%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '133.86'
$c.ClearFormats()
$ws.Range('E6').Value = '  -0.77%  '
$ws.Range('E7').Value = '  +0.54%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.561'
$c.ClearFormats()
$ws.Range('E8').Value = '  +4.46%  '
$ws.Range('E9').Value = '  +0.34%  '
$ws.Range('E10').Value = '  +1.97%  '
$ws.Range('E11').Value = '  -1.85%  '
$ws.Range('E12').Value = '  +0.32%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '23.76'
$c.ClearFormats()
$ws.Range('E13').Value = '  +0.99%  '
$ws.Range('D14').Value = '2.760.36'
$ws.Range('E14').Value = '  -0.15%  '
$ws.Range('D15').Value = '57.897.95'
$ws.Range('E15').Value = '  +0.32%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '0.0000134'
$c.ClearFormats()
$ws.Range('E16').Value = '  +0.31%  '
$ws.Range('D17').Value = '2.341.33'
$ws.Range('E17').Value = '  -0.30%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '10.65'
$c.ClearFormats()
$ws.Range('E18').Value = '  +0.83%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '4.29'
$c.ClearFormats()
$ws.Range('E19').Value = '  +1.81%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '328.32'
$c.ClearFormats()
$ws.Range('E20').Value = '  -1.95%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '6.74'
$c.ClearFormats()
$ws.Range('E21').Value = '  -0.13%  '
$ws.Range('E22').Value = '  +0.04%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '62.90'
$c.ClearFormats()
$ws.Range('E23').Value = '  +1.78%  '
$ws.Range('E24').Value = '  -3.35%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '0.995'
$c.ClearFormats()
$ws.Range('E25').Value = '  -0.37%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '8.29'
$c.ClearFormats()
$ws.Range('E26').Value = '  -1.79%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '1.33'
$c.ClearFormats()
$ws.Range('E27').Value = '  -6.48%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '1.76'
$c.ClearFormats()
$ws.Range('E28').Value = '  +0.30%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '170.27'
$c.ClearFormats()
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('D30').Value = '0.0₃0733'
$ws.Range('E30').Value = '  -0.51%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '6.12'
$c.ClearFormats()
$ws.Range('E31').Value = '  -0.86%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '18.31'
$c.ClearFormats()
$ws.Range('E32').Value = '  -1.37%  '
$ws.Range('E33').Value = '  -1.34%  '
$ws.Range('E34').Value = '  +0.03%  '
$ws.Range('E35').Value = '  +0.70%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '4.16'
$c.ClearFormats()
$ws.Range('E36').Value = '  +0.27%  '
$ws.Range('E37').Value = '  -2.47%  '
$ws.Range('E38').Value = '  -0.90%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '39.06'
$c.ClearFormats()
$ws.Range('E39').Value = '  -0.98%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '141.28'
$c.ClearFormats()
$ws.Range('E40').Value = '  -5.81%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.376'
$c.ClearFormats()
$ws.Range('E41').Value = '  -0.70%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '287.74'
$c.ClearFormats()
$ws.Range('E42').Value = '  +1.19%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '3.62'
$c.ClearFormats()
$ws.Range('E43').Value = '  +0.11%  '
$ws.Range('E45').Value = '  +0.71%  '
$ws.Range('E46').Value = '  -1.72%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.565'
$c.ClearFormats()
$ws.Range('E47').Value = '  +0.35%  '
$ws.Range('E48').Value = '  +1.40%  '
$ws.Range('E49').Value = '  +0.19%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '11.08'
$c.ClearFormats()
$ws.Range('E50').Value = '  +0.24%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '0.953'
$c.ClearFormats()
